$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Reorder the config columns C/D/E (DestinationEmail, SourceEmailCredential,
# RetryLimit) and add a brand new ScrapeAmount column (F) with its value.
# ---------------------------------------------------------------------------

# Row 1 - headers
$ws.Range("A1").Value = "Entries"
$ws.Range("B1").Value = "Description"
$ws.Range("C1").Value = "DestinationEmail"
$ws.Range("D1").Value = "SourceEmailCredential"
$ws.Range("E1").Value = "RetryLimit"
$ws.Range("F1").Value = "ScrapeAmount"

# Row 2 - values
$ws.Range("A2").Value = 400
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "burcea.bogdan.madalin@gmail.com"
$ws.Range("D2").Value = "bot_gmail"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 100

# Move the "Hyperlink" cell style off the old D2/D5 cells and onto the new
# C2/C5 cells (the mail-credential column moved from D to C).
$ws.Range("D2").Style = "Normal"

# Move the mailto hyperlink itself from D2 to C2. Re-apply the built-in
# "Hyperlink" style afterwards since Hyperlinks.Add stamps its own
# formatting on the anchor cell.
$ws.Range("D2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:burcea.bogdan.madalin@gmail.com") | Out-Null
$ws.Range("C2").Style = "Hyperlink"
$ws.Range("C5").Style = "Hyperlink"

# Clean up the now-empty D5 placeholder cell (its formatting moved to C5).
$ws.Range("D5").Clear() | Out-Null

# ---------------------------------------------------------------------------
# Column widths - values below are the closest COM-addressable ColumnWidth
# inputs that reproduce the target stored widths through Excel's whole-pixel
# rounding of the ColumnWidth property.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 15
$ws.Columns.Item(3).ColumnWidth = 38.666666666666664
$ws.Columns.Item(4).ColumnWidth = 41.166666666666664
$ws.Columns.Item(5).ColumnWidth = 19
$ws.Columns.Item(6).ColumnWidth = 20

# Update the active selection shown when the sheet was last saved.
$ws.Range("F8").Select() | Out-Null
